# Regenerate the "K" column (column G) values to reflect actual strikeout
# counts (K) instead of the previous total-strikes (Strike#) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeouts) values for rows 2-40, replacing the old Strike# values.
$newK = @{
    2  = 5
    3  = 8
    4  = 6
    5  = 6
    6  = 7
    7  = 7
    8  = 10
    9  = 4
    10 = 9
    11 = 8
    12 = 7
    13 = 6
    14 = 7
    15 = 5
    16 = 7
    17 = 8
    18 = 6
    19 = 11
    20 = 7
    21 = 4
    22 = 9
    23 = 6
    24 = 7
    25 = 12
    26 = 7
    27 = 3
    28 = 8
    29 = 9
    30 = 4
    31 = 4
    32 = 8
    33 = 5
    34 = 7
    35 = 9
    36 = 5
    37 = 6
    38 = 2
    39 = 3
    40 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
